$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# New header-row values F2/G2 go in first, so their shared strings are
# appended to the table before the A1 text edit below - this reproduces the
# same shared-string ordering as the target workbook.
$ws2.Range("F2").Value = "เป็นค่าว่างได้"
$ws2.Range("G2").Value = "เพิ่มข้อมูล ต้องไม่เป็นค่าว่าง"

# F2/G2 reuse the same look as the other row-2 header cells (A2.."Neutral"-ish
# style, B2.."Bad"-ish style) - copy the formatting only, so the existing
# cellXf entries are reused instead of new ones being created.
$ws2.Range("A2").Copy()
$ws2.Range("F2").PasteSpecial(-4122)
$ws2.Range("B2").Copy()
$ws2.Range("G2").PasteSpecial(-4122)

# A1's label picks up a clarifying suffix.
$ws2.Range("A1").Value = "เลขที่  (null)  ตัวเลขเท่านั้น"

# New merged header cell spanning F1:G1 (centered, no border/fill - matches
# the newly introduced cellXf). Style the anchor cell before merging so both
# F1 and G1 end up sharing the new centered cellXf.
$ws2.Range("F1").HorizontalAlignment = -4108
$ws2.Range("F1:G1").Merge()

# Match the width of the two new columns.
$ws2.Range("F1:G1").ColumnWidth = 25.7265625

# Sheet2 becomes the active/selected sheet (was Sheet1 before).
$ws2.Activate()
